# Natmi following Dr Hou advice
# The underlying LR-pair computation was re-run: instead of pairing each
# sending cluster with every OTHER cluster, cells are now paired with
# every cluster including themselves, which adds the missing ECs->ECs /
# FAPs->FAPs / sCs->sCs (and ECs as a sending cluster) rows, growing the
# table from a 3x2 "sending x target" grid (6 rows) to a full 3x3 grid
# (9 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Each inner array holds, in column order A..T:
#   Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
#   then the 16 numeric metric columns (E..T).
$rows = @(
    @("ECs","Adam12","Itgb1","ECs",  2,0.6666666666666666,1.958141333333333,5.874423999999999,0.1445807708852573,0.1445807708852573,3,1,112.513392,337.540176,0.3275312977368564,0.3275312977368564,220.317123428736,1982.854110858624,0.04735472751584341,0.04735472751584342),
    @("ECs","Adam12","Itgb1","FAPs", 2,0.6666666666666666,1.958141333333333,5.874423999999999,0.1445807708852573,0.1445807708852573,3,1,106.314466,318.943398,0.3094859589441663,0.3094859589441664,208.1787502058613,1873.608751852752,0.04474571852231064,0.04474571852231066),
    @("ECs","Adam12","Itgb1","sCs",  2,0.6666666666666666,1.958141333333333,5.874423999999999,0.1445807708852573,0.1445807708852573,3,1,124.6916553333333,374.074966,0.3629827433189773,0.3629827433189773,244.1638842299538,2197.474958069584,0.05248032484710319,0.05248032484710322),
    @("FAPs","Adam12","Itgb1","ECs", 3,1,5.833003000000001,17.499009,0.4306839633891008,0.4306839633891009,3,1,112.513392,337.540176,0.3275312977368564,0.3275312977368564,656.290953076176,5906.618577685584,0.1410624774432849,0.1410624774432849),
    @("FAPs","Adam12","Itgb1","FAPs",3,1,5.833003000000001,17.499009,0.4306839633891008,0.4306839633891009,3,1,106.314466,318.943398,0.3094859589441663,0.3094859589441664,620.132599121398,5581.193392092582,0.1332906394113501,0.1332906394113501),
    @("FAPs","Adam12","Itgb1","sCs", 3,1,5.833003000000001,17.499009,0.4306839633891008,0.4306839633891009,3,1,124.6916553333333,374.074966,0.3629827433189773,0.3629827433189773,727.3267996342995,6545.941196708694,0.1563308465344658,0.1563308465344658),
    @("sCs","Adam12","Itgb1","ECs",  3,1,5.752436333333333,17.257309,0.4247352657256419,0.4247352657256419,3,1,112.513392,337.540176,0.3275312977368564,0.3275312977368564,647.2261241273759,5825.035117146384,0.139114092777728,0.139114092777728),
    @("sCs","Adam12","Itgb1","FAPs", 3,1,5.752436333333333,17.257309,0.4247352657256419,0.4247352657256419,3,1,106.314466,318.943398,0.3094859589441663,0.3094859589441664,611.5671969773313,5504.104772795982,0.1314496010105056,0.1314496010105056),
    @("sCs","Adam12","Itgb1","sCs", 3,1,5.752436333333333,17.257309,0.4247352657256419,0.4247352657256419,3,1,124.6916553333333,374.074966,0.3629827433189773,0.3629827433189773,717.2808086029438,6455.527277426494,0.1541715719374083,0.1541715719374083)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowData = $rows[$r]
    $excelRow = $r + 2   # row 1 is the header
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Range($cols[$c] + $excelRow).Value = $rowData[$c]
    }
}
